# The "SASTH46@" value (used as a cohort/password token for students in
# rows 2-247) gets corrected to "SASTH46" by stripping the stray trailing
# "@". All 246 cells in B2:B247 share the same text, so writing the new
# value across the whole range at once keeps them on a single shared
# string (matching how the workbook was re-saved).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:B247")
$rng.Value = "SASTH46"

# The edit also dropped the cell formatting that had been applied to this
# column (a fill-only style override), reverting these cells to the
# worksheet's default style.
$rng.ClearFormats()

# Leave the selection where the edit was made.
$rng.Select() | Out-Null
